$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 1919.3889
$ws.Range("I11").Value = 1919.3889
$ws.Range("K11").Value = 1919.3889
$ws.Range("M11").Value = -1779.3889

$ws.Range("H62").Value = 41670340
$ws.Range("I62").Value = 66669820
$ws.Range("K62").Value = 66669820
$ws.Range("M62").Value = -66669196

$ws.Range("H65").Value = 41670340
$ws.Range("I65").Value = 66669820
$ws.Range("K65").Value = 333349100
$ws.Range("M65").Value = -333345980

$ws.Range("H87").Value = 13524.686
$ws.Range("J87").Value = 13524.686
$ws.Range("L87").Value = 13524.686
$ws.Range("N87").Value = -16020.686

$ws.Range("H90").Value = 13524.686
$ws.Range("J90").Value = 13524.686
$ws.Range("L90").Value = 40574.058
$ws.Range("N90").Value = -53054.058

$ws.Range("H106").Value = 2101.4614
$ws.Range("I106").Value = 1665.3636
$ws.Range("K106").Value = 1665.3636
$ws.Range("M106").Value = -1034.3636

$ws.Range("H113").Value = 13590.909
$ws.Range("I113").Value = 2750
$ws.Range("J113").Value = 19785.715
$ws.Range("K113").Value = 2750
$ws.Range("L113").Value = 19785.715
$ws.Range("M113").Value = 504
$ws.Range("N113").Value = -26293.715

$ws.Range("H129").Value = 740.9211
$ws.Range("I129").Value = 421.3846
$ws.Range("J129").Value = 907.08
$ws.Range("K129").Value = 1264.1538
$ws.Range("L129").Value = 2721.24
$ws.Range("M129").Value = 3735.8462
$ws.Range("N129").Value = -12721.24

$ws.Range("H132").Value = 2000.5938
$ws.Range("I132").Value = 2067.3
$ws.Range("J132").Value = 1000
$ws.Range("K132").Value = 6201.900000000001
$ws.Range("L132").Value = 3000
$ws.Range("M132").Value = -3671.900000000001
$ws.Range("N132").Value = -8060

$ws.Range("H135").Value = 1520.4445
$ws.Range("I135").Value = 516
$ws.Range("K135").Value = 4644
$ws.Range("M135").Value = -2109

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4016.6956
$ws.Range("I32").Value = 2526.6453
$ws.Range("J32").Value = 17214.285
$ws.Range("K32").Value = 2526.6453
$ws.Range("L32").Value = 17214.285
$ws.Range("M32").Value = -2239.6453
$ws.Range("N32").Value = -17788.285

$ws.Range("H45").Value = 5045
$ws.Range("I45").Value = 1090
$ws.Range("J45").Value = 9000
$ws.Range("K45").Value = 1090
$ws.Range("L45").Value = 9000
$ws.Range("M45").Value = -713
$ws.Range("N45").Value = -9754

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H21").Value = 17644.8
$ws.Range("J21").Value = 17644.8
$ws.Range("L21").Value = 17644.8
$ws.Range("N21").Value = -18116.8

$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").ClearContents()
$ws.Range("N125").Value = 0

$ws.Range("H126").Value = 52780
$ws.Range("J126").Value = 52780
$ws.Range("L126").Value = 52780
$ws.Range("N126").Value = -62660

$ws.Range("H134").Value = 2719.05
$ws.Range("I134").Value = 1708.3572
$ws.Range("K134").Value = 5125.071599999999
$ws.Range("M134").Value = -2590.071599999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4132.0884
$ws.Range("I31").Value = 929
$ws.Range("J31").Value = 6979.278
$ws.Range("K31").Value = 929
$ws.Range("L31").Value = 6979.278
$ws.Range("M31").Value = -634
$ws.Range("N31").Value = -7569.278

$ws.Range("H34").Value = 4132.0884
$ws.Range("I34").Value = 929
$ws.Range("J34").Value = 6979.278
$ws.Range("K34").Value = 929
$ws.Range("L34").Value = 6979.278
$ws.Range("M34").Value = -727
$ws.Range("N34").Value = -7383.278

$ws.Range("H99").Value = 3568616.2
$ws.Range("I99").Value = 5335924.5
$ws.Range("K99").Value = 5335924.5
$ws.Range("M99").Value = -5334426.5

$ws.Range("H126").Value = 3568616.2
$ws.Range("I126").Value = 5335924.5
$ws.Range("K126").Value = 16007773.5
$ws.Range("M126").Value = -16005303.5

$ws.Range("H134").Value = 5234.161
$ws.Range("I134").Value = 5555.7144
$ws.Range("J134").Value = 2233
$ws.Range("K134").Value = 16667.1432
$ws.Range("L134").Value = 6699
$ws.Range("M134").Value = -14132.1432
$ws.Range("N134").Value = -11769

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 549140.7
$ws.Range("I68").Value = 1344056.9
$ws.Range("J68").Value = 1531.8
$ws.Range("K68").Value = 4032170.7
$ws.Range("L68").Value = 4595.4
$ws.Range("M68").Value = -4031359.7
$ws.Range("N68").Value = -6217.4

$ws.Range("H71").Value = 549140.7
$ws.Range("I71").Value = 1344056.9
$ws.Range("J71").Value = 1531.8
$ws.Range("K71").Value = 12096512.1
$ws.Range("L71").Value = 13786.2
$ws.Range("M71").Value = -12092456.1
$ws.Range("N71").Value = -21898.2

$ws.Range("H98").Value = 538
$ws.Range("I98").Value = 615.5
$ws.Range("J98").Value = 434.66666
$ws.Range("K98").Value = 1846.5
$ws.Range("L98").Value = 1303.99998
$ws.Range("M98").Value = -348.5
$ws.Range("N98").Value = -4299.999980000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H51").Value = 0
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("L51").ClearContents()
$ws.Range("M51").ClearContents()
$ws.Range("N51").Value = 0

$ws.Range("H80").Value = 5249.9375
$ws.Range("I80").Value = 5777.778
$ws.Range("J80").Value = 4571.2856
$ws.Range("K80").Value = 5777.778
$ws.Range("L80").Value = 4571.2856
$ws.Range("M80").Value = -4779.778
$ws.Range("N80").Value = -6567.2856

$ws.Range("H83").Value = 5249.9375
$ws.Range("I83").Value = 5777.778
$ws.Range("J83").Value = 4571.2856
$ws.Range("K83").Value = 28888.89
$ws.Range("L83").Value = 22856.428
$ws.Range("M83").Value = -23896.89
$ws.Range("N83").Value = -32840.428

$ws.Range("H102").Value = 1005.5
$ws.Range("I102").Value = 1005.5
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 1005.5
$ws.Range("L102").Value = 0
$ws.Range("M102").ClearContents()
$ws.Range("N102").Value = 616.5

$ws.Range("H132").Value = 3291.48
$ws.Range("I132").Value = 2361.3635
$ws.Range("J132").Value = 4022.2856
$ws.Range("K132").Value = 7084.0905
$ws.Range("L132").Value = 12066.8568
$ws.Range("M132").Value = -4554.0905
$ws.Range("N132").Value = -17126.8568

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 22224544
$ws.Range("I122").Value = 22224544
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 66673632
$ws.Range("L122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -66671182

$ws.Range("H127").Value = 77777
$ws.Range("J127").Value = 77777
$ws.Range("L127").Value = 77777
$ws.Range("N127").Value = -87697

$ws.Range("H132").Value = 4378.5435
$ws.Range("I132").Value = 4701.28
$ws.Range("J132").Value = 3994.3333
$ws.Range("K132").Value = 14103.84
$ws.Range("L132").Value = 11982.9999
$ws.Range("M132").Value = -11573.84
$ws.Range("N132").Value = -17042.9999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 806.8182
$ws.Range("I81").Value = 625
$ws.Range("J81").Value = 910.7143
$ws.Range("K81").Value = 1250
$ws.Range("L81").Value = 1821.4286
$ws.Range("M81").Value = -189
$ws.Range("N81").Value = -3943.4286

$ws.Range("H84").Value = 806.8182
$ws.Range("I84").Value = 625
$ws.Range("J84").Value = 910.7143
$ws.Range("K84").Value = 6250
$ws.Range("L84").Value = 9107.143
$ws.Range("M84").Value = -946
$ws.Range("N84").Value = -19715.143

$ws.Range("H122").Value = 54843.95
$ws.Range("I122").Value = 85345
$ws.Range("K122").Value = 256035
$ws.Range("M122").Value = -253585

$ws.Range("H132").Value = 1843.1321
$ws.Range("I132").Value = 1349.4
$ws.Range("J132").Value = 2487.1304
$ws.Range("K132").Value = 4048.2
$ws.Range("L132").Value = 7461.3912
$ws.Range("M132").Value = -1518.2
$ws.Range("N132").Value = -12521.3912

$ws.Range("H136").Value = 4507409.5
$ws.Range("I136").Value = 7576286.5
$ws.Range("J136").Value = 6390.1665
$ws.Range("K136").Value = 22728859.5
$ws.Range("L136").Value = 19170.4995
$ws.Range("M136").Value = -22726309.5

